$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (D), Volume% (E) and Hora (G) columns hold text that looks numeric.
# A leading apostrophe forces Excel to keep the literal text (e.g. "20",
# "2.91%", "0.0002000") instead of auto-converting it to a number/percent.

$ws.Range('D2').Value = '''297.82'
$ws.Range('E2').Value = '''2.91%'
$ws.Range('G2').Value = '''20'

$ws.Range('D3').Value = '''41.52'
$ws.Range('E3').Value = '''3.15%'
$ws.Range('G3').Value = '''20'

$ws.Range('D4').Value = '''5.027'
$ws.Range('E4').Value = '''-0.45%'
$ws.Range('G4').Value = '''20'

$ws.Range('D5').Value = '''0.07542'
$ws.Range('E5').Value = '''3.41%'
$ws.Range('G5').Value = '''20'

$ws.Range('B6').Value = 'FTXToken'
$ws.Range('C6').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D6').Value = '''1.575'
$ws.Range('E6').Value = '''2.08%'
$ws.Range('G6').Value = '''20'

$ws.Range('B7').Value = 'MXToken'
$ws.Range('C7').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D7').Value = '''0.9311'
$ws.Range('E7').Value = '''1.39%'
$ws.Range('G7').Value = '''20'

$ws.Range('B8').Value = 'BTSEToken'
$ws.Range('C8').Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range('D8').Value = '''2.425'
$ws.Range('E8').Value = '''1.17%'
$ws.Range('G8').Value = '''20'

$ws.Range('B9').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C9').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('D9').Value = '''0.1191'
$ws.Range('E9').Value = '''1.20%'
$ws.Range('G9').Value = '''20'

$ws.Range('B10').Value = 'WazirX'
$ws.Range('C10').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('D10').Value = '''0.1844'
$ws.Range('E10').Value = '''7.22%'
$ws.Range('G10').Value = '''20'

$ws.Range('B11').Value = 'MandalaExchangeToken'
$ws.Range('C11').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('D11').Value = '''0.08914'
$ws.Range('E11').Value = '''3.31%'
$ws.Range('G11').Value = '''20'

$ws.Range('B12').Value = 'BitrueCoin'
$ws.Range('C12').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('D12').Value = '''0.04050'
$ws.Range('E12').Value = '''-3.08%'
$ws.Range('G12').Value = '''20'

$ws.Range('B13').Value = 'BitMartToken'
$ws.Range('C13').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D13').Value = '''0.1056'
$ws.Range('E13').Value = '''0.16%'
$ws.Range('G13').Value = '''20'

$ws.Range('B14').Value = 'BitForexToken'
$ws.Range('C14').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D14').Value = '''0.001281'
$ws.Range('E14').Value = '''1.17%'
$ws.Range('G14').Value = '''20'

$ws.Range('B15').Value = 'TigerCash'
$ws.Range('C15').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('D15').Value = '''0.005881'
$ws.Range('E15').Value = '''1.56%'
$ws.Range('G15').Value = '''20'

$ws.Range('B16').Value = 'LEO'
$ws.Range('C16').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D16').Value = '''3.344'
$ws.Range('E16').Value = '''-1.41%'
$ws.Range('G16').Value = '''20'

$ws.Range('B17').Value = 'GateToken'
$ws.Range('C17').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range('D17').Value = '''4.376'
$ws.Range('E17').Value = '''2.16%'
$ws.Range('G17').Value = '''20'

$ws.Range('D18').Value = '''0.3311'
$ws.Range('E18').Value = '''-0.21%'
$ws.Range('G18').Value = '''20'

$ws.Range('D19').Value = '''7.894'
$ws.Range('E19').Value = '''0.52%'
$ws.Range('G19').Value = '''20'

$ws.Range('D20').Value = '''0.1419'
$ws.Range('E20').Value = '''4.87%'
$ws.Range('G20').Value = '''20'

$ws.Range('E21').Value = '''3.74%'
$ws.Range('G21').Value = '''20'

$ws.Range('D22').Value = '''0.04046'
$ws.Range('E22').Value = '''4.79%'
$ws.Range('G22').Value = '''20'

$ws.Range('D23').Value = '''0.001264'
$ws.Range('E23').Value = '''-0.39%'
$ws.Range('G23').Value = '''20'

$ws.Range('D24').Value = '''0.004172'
$ws.Range('E24').Value = '''9.03%'
$ws.Range('G24').Value = '''20'

$ws.Range('E25').Value = '''-3.97%'
$ws.Range('G25').Value = '''20'

$ws.Range('E26').Value = '''-0.15%'
$ws.Range('G26').Value = '''20'

$ws.Range('G27').Value = '''20'

$ws.Range('G28').Value = '''20'

$ws.Range('G29').Value = '''20'

$ws.Range('G30').Value = '''20'

$ws.Range('G31').Value = '''20'

$ws.Range('G32').Value = '''20'

$ws.Range('G33').Value = '''20'

$ws.Range('G34').Value = '''20'

$ws.Range('G35').Value = '''20'

$ws.Range('G36').Value = '''20'

$ws.Range('G37').Value = '''20'

$ws.Range('D38').Value = '''0.02416'
$ws.Range('E38').Value = '''4.56%'
$ws.Range('G38').Value = '''20'

$ws.Range('D39').Value = '''0.05221'
$ws.Range('E39').Value = '''4.91%'
$ws.Range('G39').Value = '''20'

$ws.Range('D40').Value = '''0.006216'
$ws.Range('E40').Value = '''-6.24%'
$ws.Range('G40').Value = '''20'

$ws.Range('D41').Value = '''0.007801'
$ws.Range('E41').Value = '''1.75%'
$ws.Range('G41').Value = '''20'

$ws.Range('E42').Value = '''4.55%'
$ws.Range('G42').Value = '''20'

$ws.Range('D43').Value = '''0.007354'
$ws.Range('E43').Value = '''-0.31%'
$ws.Range('G43').Value = '''20'

$ws.Range('D44').Value = '''0.007828'
$ws.Range('E44').Value = '''11.09%'
$ws.Range('G44').Value = '''20'

$ws.Range('D45').Value = '''0.2982'
$ws.Range('E45').Value = '''-4.43%'
$ws.Range('G45').Value = '''20'

$ws.Range('D46').Value = '''0.00006372'
$ws.Range('E46').Value = '''-1.15%'
$ws.Range('G46').Value = '''20'

$ws.Range('D47').Value = '''0.00000000750'
$ws.Range('E47').Value = '''-0.22%'
$ws.Range('G47').Value = '''20'

$ws.Range('B48').Value = 'CoinbaseStockToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin'
$ws.Range('D48').Value = '''0.004199'
$ws.Range('E48').Value = '''-0.05%'
$ws.Range('G48').Value = '''20'

$ws.Range('B49').Value = 'BOLO'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ogrGe0dEab+bolo-bolo'
$ws.Range('D49').Value = '''0.04431'
$ws.Range('E49').Value = '''418.99%'
$ws.Range('G49').Value = '''20'

$ws.Range('D50').Value = '''0.00002100'
$ws.Range('E50').Value = '''-0.22%'
$ws.Range('G50').Value = '''20'

$ws.Range('D51').Value = '''0.0002000'
$ws.Range('E51').Value = '''-0.22%'
$ws.Range('G51').Value = '''20'
